# Apply the "BFO to IES mapping" update: replace the single sample row
# (BFO_0000027 / DispositionalClass, with hyperlinks) with a richer row
# describing BFO_0000015 (process) <-> ies:Event, and give the sheet its
# full column layout/widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks (and their formatting baggage) -------------
$ws.Hyperlinks.Delete()

# --- Reset row 2 to a plain/default look before writing new values --------
$ws.Range("A2:H2").Style = "Normal"

# --- Write the new data row -------------------------------------------------
$ws.Range("A2").Value2 = "http://purl.obolibrary.org/obo/BFO_0000015"
$ws.Range("B2").Value2 = "process"
$ws.Range("C2").Value2 = "R:only×3"
$ws.Range("D2").Value2 = "SubClassOf: ns1:BFO_0000003 | SubClassOf: ns1:BFO_0000117 only (ns1:BFO_0000015 OR ns1:BFO_0000035) | SubClassOf: ns1:BFO_0000132 only ns1:BFO_0000015 | SubClassOf: ns1:BFO_0000139 only ns1:BFO_0000015"
$ws.Range("E2").Value2 = "http://ies.data.gov.uk/ontology/ies4#Event"
$ws.Range("F2").Value2 = "Event"
$ws.Range("G2").Value2 = "R:only"
$ws.Range("H2").Value2 = "SubClassOf: ies:Element | SubClassOf: ns1:BFO_0000015 | SubClassOf: ns1:BFO_0000178 only (ns1:BFO_0000029 OR ns1:BFO_0000140)"

# --- Row 2 sizing ------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15

# --- Column layout / widths --------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 72
$ws.Columns.Item(2).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 46
$ws.Columns.Item(6).ColumnWidth = 48.67
$ws.Columns.Item(8).ColumnWidth = 25.83
$ws.Columns.Item(9).ColumnWidth = 18.17

# --- Touch the sheet's last row (mirrors the source workbook's full used
#     range / default-row-height bookkeeping) --------------------------------
$ws.Rows.Item(1048576).RowHeight = 12.8

# --- Selection lands on A2, same as in the edited workbook -------------------
$ws.Range("A2").Select()
